$d = $word.ActiveDocument

# --- 1. Consolidate the three runs that make up the opening sentence of the
#        third paragraph into a single run (pure run-merge, no text change).
$d.Content.Find.Execute( `
    "This script calculates the sensitivity of the excess generation factor to relative sizing of annual generation to load. The original calculation assumes that ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "This script calculates the sensitivity of the excess generation factor to relative sizing of annual generation to load. The original calculation assumes that ", `
    2) | Out-Null

# --- 2. Merge "hourly_l" + "oad" into a single "hourly_load" run.
$d.Content.Find.Execute("hourly_load", `
    $true, $false, $false, $false, $false, $true, 1, $false, "hourly_load", 2) | Out-Null

# --- 3. Merge ") for an annual time series. " + "That value is the " into one run.
$d.Content.Find.Execute(") for an annual time series. That value is the ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ") for an annual time series. That value is the ", 2) | Out-Null

# --- 4. Merge "...value cha" with the following run so the whole sentence
#        ending in "...system is smaller or larger?" lives in one run again
#        (mirrors the pre-edit layout), then carve the real edit out of it:
#        the typo fix is a single character substitution, "is" -> "if", i.e.
#        only the "s" becomes "f". Word naturally leaves that freshly-typed
#        "f" in its own run, with the "_GoBack" bookmark marking the last
#        edit location right after it.
$d.Content.Find.Execute("value cha", $true, $false, $false, $false, $false, $true, 1, $false, "value cha", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("nge is the system is smaller or larger?", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $rng.Start            # start of "nge is the system..."
$sPos = $base + 5             # the 's' in 'is' ("n-g-e-space-i-s")

# Split the run right before the 's' (isolates it) using a temporary bookmark.
$tempSplit = $d.Range($sPos, $sPos)
$d.Bookmarks.Add("TempSplit", $tempSplit) | Out-Null

# Replace the isolated 's' with 'f' -> "is" becomes "if".
$sRange = $d.Range($sPos, $sPos + 1)
$sRange.Text = "f"

# Remove the temporary bookmark (the run split it forced remains in place).
$d.Bookmarks("TempSplit").Delete()

# --- 5. Merge the "bounded" / "by 50% to 115% ..." / "The relationship is " runs.
$d.Content.Find.Execute("bounded by 50% to 115% of annual load. The relationship is ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "bounded by 50% to 115% of annual load. The relationship is ", 2) | Out-Null

# --- 6. Move the "_GoBack" bookmark from its old spot (between "t a" and
#        " small sizing ...") to right after the newly-typed "f".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$newBm = $d.Range($sPos + 1, $sPos + 1)
$d.Bookmarks.Add("_GoBack", $newBm) | Out-Null
